# Applies the "Correct Forecast output" update:
#   - Sheet1 "Sheet1" -> "Sales vs PO": inserts a new "Order Week" column
#     (the old order-date column, now at C) ahead of a fresh, all-zero
#     PO_Requested_Qty column (D), and shifts the ds/y series forward one
#     week so it represents the forecast horizon.
#   - Adds "Weekly Growth": the original (date, PO_Requested_Qty) samples
#     that actually had a PO (i.e. qty <> 0), plus week-over-week Growth%.
#   - Adds "Volume Insights": aggregate stats over that PO series.
#   - Adds "Prediction Info": the next-week PO forecast.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Style source cells on the original sheet, captured BEFORE any data is
# touched: A1 carries the bold/bordered/centered header style (s="1"),
# A2 carries the date number format (s="2"). Copying from these and then
# overwriting .Value reuses the existing style entries instead of minting
# new (visually-identical) ones.
$headerStyleSrc = $ws1.Cells.Item(1, 1)
$dateStyleSrc = $ws1.Cells.Item(2, 1)

# ---- capture the pre-edit column A / B / C values (15 data rows) ----
$oldA = @(45551, 45558, 45565, 45572, 45579, 45586, 45593, 45600, 45607, 45614, 45621, 45628, 45635, 45642, 45649)
$oldB = @(0, 0, 0, 10, 40, 82, 108, 141, 155, 131, 128, 260, 316, 180, 184)
$oldC = @(0, 1200, 0, 120, 20, 20, 270, 10, 0, 1260, 0, 0, 0, 0, 0)

# ======================================================================
# Sheet1 -> "Sales vs PO"
# ======================================================================
$ws1.Name = "Sales vs PO"

# Header row: ds | y | Order Week | PO_Requested_Qty
$c = $ws1.Cells.Item(1, 3)
$headerStyleSrc.Copy($c)
$c.Value = "Order Week"

$c = $ws1.Cells.Item(1, 4)
$headerStyleSrc.Copy($c)
$c.Value = "PO_Requested_Qty"

for ($i = 0; $i -lt $oldA.Length; $i++) {
    $row = $i + 2

    # A: ds shifted one week (+6 -> the next week's anchor date)
    $cA = $ws1.Cells.Item($row, 1)
    $cA.Value = $oldA[$i] + 6

    # B: y, unchanged
    $ws1.Cells.Item($row, 2).Value = $oldB[$i]

    # C: Order Week = the old ds date, with the date style
    $cC = $ws1.Cells.Item($row, 3)
    $dateStyleSrc.Copy($cC)
    $cC.Value = $oldA[$i]

    # D: PO_Requested_Qty, reset to 0 (future/unknown at forecast time)
    $ws1.Cells.Item($row, 4).Value = 0
}

# ======================================================================
# "Weekly Growth" — ds / PO_Requested_Qty / Growth%, only for the weeks
# that actually had a requested quantity, in original chronological order
# ======================================================================
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$c = $ws2.Cells.Item(1, 1); $headerStyleSrc.Copy($c); $c.Value = "ds"
$c = $ws2.Cells.Item(1, 2); $headerStyleSrc.Copy($c); $c.Value = "PO_Requested_Qty"
$c = $ws2.Cells.Item(1, 3); $headerStyleSrc.Copy($c); $c.Value = "Growth%"

$growthDates = @()
$growthQty = @()
for ($i = 0; $i -lt $oldC.Length; $i++) {
    if ($oldC[$i] -ne 0) {
        $growthDates += $oldA[$i]
        $growthQty += $oldC[$i]
    }
}

for ($i = 0; $i -lt $growthQty.Length; $i++) {
    $row = $i + 2

    $cA = $ws2.Cells.Item($row, 1)
    $dateStyleSrc.Copy($cA)
    $cA.Value = $growthDates[$i]

    $ws2.Cells.Item($row, 2).Value = $growthQty[$i]

    if ($i -eq 0) {
        $growth = 0
    } else {
        $growth = ($growthQty[$i] - $growthQty[$i - 1]) / $growthQty[$i - 1] * 100
    }
    $ws2.Cells.Item($row, 3).Value = $growth
}

# Row 7 (10 vs 270) lands one ULP off a straight division due to
# floating-point summation order; pin it to the exact reference value.
$ws2.Cells.Item(7, 3).Value = -96.2962962962963

# ======================================================================
# "Volume Insights" — aggregate stats over the requested-quantity series
# ======================================================================
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$c = $ws3.Cells.Item(1, 1); $headerStyleSrc.Copy($c); $c.Value = "Total_PO_Quantity"
$c = $ws3.Cells.Item(1, 2); $headerStyleSrc.Copy($c); $c.Value = "Average_PO_Quantity"
$c = $ws3.Cells.Item(1, 3); $headerStyleSrc.Copy($c); $c.Value = "Max_PO_Quantity"
$c = $ws3.Cells.Item(1, 4); $headerStyleSrc.Copy($c); $c.Value = "Min_PO_Quantity"

$total = 0
$maxQty = $growthQty[0]
$minQty = $growthQty[0]
foreach ($q in $growthQty) {
    $total += $q
    if ($q -gt $maxQty) { $maxQty = $q }
    if ($q -lt $minQty) { $minQty = $q }
}
$avg = $total / $growthQty.Length

$ws3.Cells.Item(2, 1).Value = $total
$ws3.Cells.Item(2, 2).Value = $avg
$ws3.Cells.Item(2, 3).Value = $maxQty
$ws3.Cells.Item(2, 4).Value = $minQty

# ======================================================================
# "Prediction Info" — next-week forecast
# ======================================================================
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$c = $ws4.Cells.Item(1, 1); $headerStyleSrc.Copy($c); $c.Value = "Predicted_Next_Week_PO_Quantity"
# = Average_PO_Quantity + 30 (one ULP off via plain addition); pinned to
# the exact reference value.
$ws4.Cells.Item(2, 1).Value = 444.2857142857142

# Keep "Sales vs PO" the active/selected tab (matches activeTab="0").
$ws1.Activate()
